# Edit script: corrects mis-ordered Eredivisie 2023-2024 match rows
# and appends 4 new matches (rows 90-93) scraped 05-11-2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows where match data (columns F:V) had been written to the wrong row ---
# Swap match data between row 8 and row 9 (index column A is left untouched)
$ws.Cells.Item(8,6).Value = "AZ Alkmaar"
$ws.Cells.Item(8,7).Value = 5
$ws.Cells.Item(8,8).Value = "G.A. Eagles"
$ws.Cells.Item(8,9).Value = 1
$ws.Cells.Item(8,10).Value = 1.37
$ws.Cells.Item(8,11).Value = "07/07/2023 11:12"
$ws.Cells.Item(8,12).Value = 1.54
$ws.Cells.Item(8,13).Value = "13/08/2023 14:27"
$ws.Cells.Item(8,14).Value = 5.51
$ws.Cells.Item(8,15).Value = "07/07/2023 11:12"
$ws.Cells.Item(8,16).Value = 4.41
$ws.Cells.Item(8,17).Value = "13/08/2023 14:26"
$ws.Cells.Item(8,18).Value = 7.89
$ws.Cells.Item(8,19).Value = "07/07/2023 11:12"
$ws.Cells.Item(8,20).Value = 6.35
$ws.Cells.Item(8,21).Value = "13/08/2023 14:26"
$ws.Cells.Item(8,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/az-alkmaar-g-a-eagles/SrIpNFQ5/"
$ws.Cells.Item(9,6).Value = "Feyenoord"
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(9,8).Value = "Sittard"
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(9,10).Value = 1.22
$ws.Cells.Item(9,11).Value = "07/07/2023 11:12"
$ws.Cells.Item(9,12).Value = 1.19
$ws.Cells.Item(9,13).Value = "13/08/2023 14:25"
$ws.Cells.Item(9,14).Value = 7.61
$ws.Cells.Item(9,15).Value = "07/07/2023 11:12"
$ws.Cells.Item(9,16).Value = 7.38
$ws.Cells.Item(9,17).Value = "13/08/2023 14:29"
$ws.Cells.Item(9,18).Value = 12.86
$ws.Cells.Item(9,19).Value = "07/07/2023 11:12"
$ws.Cells.Item(9,20).Value = 16.76
$ws.Cells.Item(9,21).Value = "13/08/2023 14:29"
$ws.Cells.Item(9,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/feyenoord-sittard/6ZGlMZuC/"

# Swap match data between row 17 and row 18 (index column A is left untouched)
$ws.Cells.Item(17,6).Value = "Sparta Rotterdam"
$ws.Cells.Item(17,7).Value = 2
$ws.Cells.Item(17,8).Value = "Feyenoord"
$ws.Cells.Item(17,9).Value = 2
$ws.Cells.Item(17,10).Value = 5
$ws.Cells.Item(17,11).Value = "13/08/2023 14:42"
$ws.Cells.Item(17,12).Value = 4.94
$ws.Cells.Item(17,13).Value = "20/08/2023 14:28"
$ws.Cells.Item(17,14).Value = 4.52
$ws.Cells.Item(17,15).Value = "13/08/2023 14:42"
$ws.Cells.Item(17,16).Value = 4.35
$ws.Cells.Item(17,17).Value = "20/08/2023 14:29"
$ws.Cells.Item(17,18).Value = 1.61
$ws.Cells.Item(17,19).Value = "13/08/2023 14:42"
$ws.Cells.Item(17,20).Value = 1.67
$ws.Cells.Item(17,21).Value = "20/08/2023 14:28"
$ws.Cells.Item(17,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-feyenoord/21WSZhAB/"
$ws.Cells.Item(18,6).Value = "Twente"
$ws.Cells.Item(18,7).Value = 3
$ws.Cells.Item(18,8).Value = "Zwolle"
$ws.Cells.Item(18,9).Value = 1
$ws.Cells.Item(18,10).Value = 1.35
$ws.Cells.Item(18,11).Value = "13/08/2023 17:12"
$ws.Cells.Item(18,12).Value = 1.44
$ws.Cells.Item(18,13).Value = "20/08/2023 14:21"
$ws.Cells.Item(18,14).Value = 5.56
$ws.Cells.Item(18,15).Value = "13/08/2023 17:12"
$ws.Cells.Item(18,16).Value = 5.21
$ws.Cells.Item(18,17).Value = "20/08/2023 14:28"
$ws.Cells.Item(18,18).Value = 8.789999999999999
$ws.Cells.Item(18,19).Value = "13/08/2023 17:12"
$ws.Cells.Item(18,20).Value = 6.95
$ws.Cells.Item(18,21).Value = "20/08/2023 14:29"
$ws.Cells.Item(18,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/twente-zwolle/KxsOzZf5/"

# Swap match data between row 36 and row 37 (index column A is left untouched)
$ws.Cells.Item(36,6).Value = "PSV"
$ws.Cells.Item(36,7).Value = 4
$ws.Cells.Item(36,8).Value = "Nijmegen"
$ws.Cells.Item(36,9).Value = 0
$ws.Cells.Item(36,10).Value = 1.24
$ws.Cells.Item(36,11).Value = "04/09/2023 08:43"
$ws.Cells.Item(36,12).Value = 1.15
$ws.Cells.Item(36,13).Value = "16/09/2023 19:38"
$ws.Cells.Item(36,14).Value = 7.01
$ws.Cells.Item(36,15).Value = "04/09/2023 08:43"
$ws.Cells.Item(36,16).Value = 9.119999999999999
$ws.Cells.Item(36,17).Value = "16/09/2023 19:57"
$ws.Cells.Item(36,18).Value = 11.5
$ws.Cells.Item(36,19).Value = "04/09/2023 08:43"
$ws.Cells.Item(36,20).Value = 17.78
$ws.Cells.Item(36,21).Value = "16/09/2023 19:57"
$ws.Cells.Item(36,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/psv-nijmegen/8zF9rNhq/"
$ws.Cells.Item(37,6).Value = "Sittard"
$ws.Cells.Item(37,7).Value = 3
$ws.Cells.Item(37,8).Value = "FC Volendam"
$ws.Cells.Item(37,9).Value = 1
$ws.Cells.Item(37,10).Value = 1.67
$ws.Cells.Item(37,11).Value = "04/09/2023 08:43"
$ws.Cells.Item(37,12).Value = 1.47
$ws.Cells.Item(37,13).Value = "16/09/2023 19:59"
$ws.Cells.Item(37,14).Value = 4.44
$ws.Cells.Item(37,15).Value = "04/09/2023 08:43"
$ws.Cells.Item(37,16).Value = 4.96
$ws.Cells.Item(37,17).Value = "16/09/2023 19:59"
$ws.Cells.Item(37,18).Value = 4.56
$ws.Cells.Item(37,19).Value = "04/09/2023 08:43"
$ws.Cells.Item(37,20).Value = 6.75
$ws.Cells.Item(37,21).Value = "16/09/2023 19:59"
$ws.Cells.Item(37,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/sittard-fc-volendam/StXtm3pS/"

# Swap match data between row 49 and row 50 (index column A is left untouched)
$ws.Cells.Item(49,6).Value = "Zwolle"
$ws.Cells.Item(49,7).Value = 0
$ws.Cells.Item(49,8).Value = "AZ Alkmaar"
$ws.Cells.Item(49,9).Value = 3
$ws.Cells.Item(49,10).Value = 4.72
$ws.Cells.Item(49,11).Value = "17/09/2023 16:13"
$ws.Cells.Item(49,12).Value = 4.63
$ws.Cells.Item(49,13).Value = "24/09/2023 16:39"
$ws.Cells.Item(49,14).Value = 4.32
$ws.Cells.Item(49,15).Value = "17/09/2023 16:13"
$ws.Cells.Item(49,16).Value = 4.03
$ws.Cells.Item(49,17).Value = "24/09/2023 16:39"
$ws.Cells.Item(49,18).Value = 1.67
$ws.Cells.Item(49,19).Value = "17/09/2023 16:13"
$ws.Cells.Item(49,20).Value = 1.77
$ws.Cells.Item(49,21).Value = "24/09/2023 16:39"
$ws.Cells.Item(49,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/zwolle-az-alkmaar/ATQ3hbM7/"
$ws.Cells.Item(50,6).Value = "Waalwijk"
$ws.Cells.Item(50,7).Value = 1
$ws.Cells.Item(50,8).Value = "Twente"
$ws.Cells.Item(50,9).Value = 0
$ws.Cells.Item(50,10).Value = 4.49
$ws.Cells.Item(50,11).Value = "17/09/2023 13:43"
$ws.Cells.Item(50,12).Value = 6.03
$ws.Cells.Item(50,13).Value = "24/09/2023 16:43"
$ws.Cells.Item(50,14).Value = 4.43
$ws.Cells.Item(50,15).Value = "17/09/2023 13:43"
$ws.Cells.Item(50,16).Value = 4.72
$ws.Cells.Item(50,17).Value = "24/09/2023 16:44"
$ws.Cells.Item(50,18).Value = 1.68
$ws.Cells.Item(50,19).Value = "17/09/2023 13:43"
$ws.Cells.Item(50,20).Value = 1.53
$ws.Cells.Item(50,21).Value = "24/09/2023 16:38"
$ws.Cells.Item(50,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-twente/OWEoHsa8/"

# Swap match data between row 69 and row 70 (index column A is left untouched)
$ws.Cells.Item(69,6).Value = "Ajax"
$ws.Cells.Item(69,7).Value = 1
$ws.Cells.Item(69,8).Value = "AZ Alkmaar"
$ws.Cells.Item(69,9).Value = 2
$ws.Cells.Item(69,10).Value = 2.24
$ws.Cells.Item(69,11).Value = "01/10/2023 16:12"
$ws.Cells.Item(69,12).Value = 2.7
$ws.Cells.Item(69,13).Value = "08/10/2023 14:29"
$ws.Cells.Item(69,14).Value = 3.77
$ws.Cells.Item(69,15).Value = "01/10/2023 16:12"
$ws.Cells.Item(69,16).Value = 3.74
$ws.Cells.Item(69,17).Value = "08/10/2023 14:29"
$ws.Cells.Item(69,18).Value = 3.06
$ws.Cells.Item(69,19).Value = "01/10/2023 16:12"
$ws.Cells.Item(69,20).Value = 2.56
$ws.Cells.Item(69,21).Value = "08/10/2023 14:29"
$ws.Cells.Item(69,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/ajax-az-alkmaar/jHLM4SBU/"
$ws.Cells.Item(70,6).Value = "Sittard"
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = "Twente"
$ws.Cells.Item(70,9).Value = 3
$ws.Cells.Item(70,10).Value = 4.64
$ws.Cells.Item(70,11).Value = "01/10/2023 16:12"
$ws.Cells.Item(70,12).Value = 4.8
$ws.Cells.Item(70,13).Value = "08/10/2023 14:29"
$ws.Cells.Item(70,14).Value = 4.14
$ws.Cells.Item(70,15).Value = "01/10/2023 16:12"
$ws.Cells.Item(70,16).Value = 4.04
$ws.Cells.Item(70,17).Value = "08/10/2023 14:29"
$ws.Cells.Item(70,18).Value = 1.71
$ws.Cells.Item(70,19).Value = "01/10/2023 16:12"
$ws.Cells.Item(70,20).Value = 1.74
$ws.Cells.Item(70,21).Value = "08/10/2023 14:29"
$ws.Cells.Item(70,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/sittard-twente/SxcSPBst/"

# --- Append 4 newly scraped matches as rows 90-93 ---
# Clone formatting (bold/border style for the index col, date format for col E) from the last existing row
$ws.Range("A89:V89").Copy()
$ws.Range("A90:V93").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 90
$ws.Cells.Item(90,1).Value = 89
$ws.Cells.Item(90,2).Value = "netherlands"
$ws.Cells.Item(90,3).Value = "eredivisie"
$ws.Cells.Item(90,4).Value = "2023-2024"
$ws.Cells.Item(90,5).Value = 45234.6875
$ws.Cells.Item(90,6).Value = "Heracles"
$ws.Cells.Item(90,7).Value = 0
$ws.Cells.Item(90,8).Value = "PSV"
$ws.Cells.Item(90,9).Value = 6
$ws.Cells.Item(90,10).Value = 7.64
$ws.Cells.Item(90,11).Value = "29/10/2023 14:42"
$ws.Cells.Item(90,12).Value = 11.48
$ws.Cells.Item(90,13).Value = "04/11/2023 16:21"
$ws.Cells.Item(90,14).Value = 6
$ws.Cells.Item(90,15).Value = "29/10/2023 14:42"
$ws.Cells.Item(90,16).Value = 7.91
$ws.Cells.Item(90,17).Value = "04/11/2023 16:21"
$ws.Cells.Item(90,18).Value = 1.34
$ws.Cells.Item(90,19).Value = "29/10/2023 14:42"
$ws.Cells.Item(90,20).Value = 1.21
$ws.Cells.Item(90,21).Value = "04/11/2023 16:12"
$ws.Cells.Item(90,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/heracles-psv/QBs4nQfA/"

# Row 91
$ws.Cells.Item(91,1).Value = 90
$ws.Cells.Item(91,2).Value = "netherlands"
$ws.Cells.Item(91,3).Value = "eredivisie"
$ws.Cells.Item(91,4).Value = "2023-2024"
$ws.Cells.Item(91,5).Value = 45234.78125
$ws.Cells.Item(91,6).Value = "Excelsior"
$ws.Cells.Item(91,7).Value = 1
$ws.Cells.Item(91,8).Value = "AZ Alkmaar"
$ws.Cells.Item(91,9).Value = 1
$ws.Cells.Item(91,10).Value = 5.11
$ws.Cells.Item(91,11).Value = "29/10/2023 17:13"
$ws.Cells.Item(91,12).Value = 9.11
$ws.Cells.Item(91,13).Value = "04/11/2023 18:43"
$ws.Cells.Item(91,14).Value = 4.63
$ws.Cells.Item(91,15).Value = "29/10/2023 17:13"
$ws.Cells.Item(91,16).Value = 5.56
$ws.Cells.Item(91,17).Value = "04/11/2023 18:43"
$ws.Cells.Item(91,18).Value = 1.6
$ws.Cells.Item(91,19).Value = "29/10/2023 17:13"
$ws.Cells.Item(91,20).Value = 1.35
$ws.Cells.Item(91,21).Value = "04/11/2023 18:35"
$ws.Cells.Item(91,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-az-alkmaar/UqlelnQc/"

# Row 92
$ws.Cells.Item(92,1).Value = 91
$ws.Cells.Item(92,2).Value = "netherlands"
$ws.Cells.Item(92,3).Value = "eredivisie"
$ws.Cells.Item(92,4).Value = "2023-2024"
$ws.Cells.Item(92,5).Value = 45234.78125
$ws.Cells.Item(92,6).Value = "Waalwijk"
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = "Feyenoord"
$ws.Cells.Item(92,9).Value = 2
$ws.Cells.Item(92,10).Value = 8.58
$ws.Cells.Item(92,11).Value = "29/10/2023 12:42"
$ws.Cells.Item(92,12).Value = 14.74
$ws.Cells.Item(92,13).Value = "04/11/2023 18:40"
$ws.Cells.Item(92,14).Value = 5.63
$ws.Cells.Item(92,15).Value = "29/10/2023 12:42"
$ws.Cells.Item(92,16).Value = 8.98
$ws.Cells.Item(92,17).Value = "04/11/2023 18:40"
$ws.Cells.Item(92,18).Value = 1.33
$ws.Cells.Item(92,19).Value = "29/10/2023 12:42"
$ws.Cells.Item(92,20).Value = 1.17
$ws.Cells.Item(92,21).Value = "04/11/2023 18:34"
$ws.Cells.Item(92,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-feyenoord/C8w0m6u4/"

# Row 93
$ws.Cells.Item(93,1).Value = 92
$ws.Cells.Item(93,2).Value = "netherlands"
$ws.Cells.Item(93,3).Value = "eredivisie"
$ws.Cells.Item(93,4).Value = "2023-2024"
$ws.Cells.Item(93,5).Value = 45234.875
$ws.Cells.Item(93,6).Value = "G.A. Eagles"
$ws.Cells.Item(93,7).Value = 5
$ws.Cells.Item(93,8).Value = "Vitesse"
$ws.Cells.Item(93,9).Value = 1
$ws.Cells.Item(93,10).Value = 1.88
$ws.Cells.Item(93,11).Value = "28/10/2023 20:13"
$ws.Cells.Item(93,12).Value = 2.02
$ws.Cells.Item(93,13).Value = "04/11/2023 20:12"
$ws.Cells.Item(93,14).Value = 3.91
$ws.Cells.Item(93,15).Value = "28/10/2023 20:13"
$ws.Cells.Item(93,16).Value = 3.8
$ws.Cells.Item(93,17).Value = "04/11/2023 20:12"
$ws.Cells.Item(93,18).Value = 4.06
$ws.Cells.Item(93,19).Value = "28/10/2023 20:13"
$ws.Cells.Item(93,20).Value = 3.69
$ws.Cells.Item(93,21).Value = "04/11/2023 20:12"
$ws.Cells.Item(93,22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/g-a-eagles-vitesse/KYt8opAG/"
